# "docs and some small updates"
#
# Target changes (per the XML diff):
#   1. "Bliss" -> "Mercy" in the "Bliss (Summons)" heading, keeping the
#      existing run split between "Mercy" and " (Summons)".
#   2. The (hidden) "_GoBack" bookmark moves from between "Increase Sum"
#      and "mon " (in the "Increase Summon Melee Damage" bullet) to sit
#      right after "Mercy" (i.e. between "Mercy" and " (Summons)").
#   3. "Increase Sum" + "mon " merge back into a single run
#      "Increase Summon ", while "Melee Damage" remains its own run.

$d = $word.ActiveDocument

# --- 1) Bliss -> Mercy ------------------------------------------------
# Scope the Find to just that one paragraph so nothing else in the
# document is touched, and only replace the first (only) match.
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Bliss*") {
        $headingPara = $p
        break
    }
}
$headingPara.Range.Find.Execute("Bliss", $true, $false, $false, $false, $false, $true, 1, $false, "Mercy", 1) | Out-Null

# --- 2) Move the _GoBack bookmark to just after "Mercy" ---------------
# Bookmarks.Add with an existing name relocates it (removing the old
# occurrence), matching Word's "names are unique" behaviour.
$full = $d.Content.Text
$mercyIdx = $full.IndexOf("Mercy")
$newBookmarkSpot = $d.Range($mercyIdx + 5, $mercyIdx + 5)
$d.Bookmarks.Add("_GoBack", $newBookmarkSpot) | Out-Null

# --- 3) Re-join "Increase Sum" + "mon " into "Increase Summon " -------
# Find the paragraph that still contains the split "Increase Sum" / "mon "
# runs (i.e. the one with "Melee Damage" after it).
$meleePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Increase Sum*Melee Damage*") {
        $meleePara = $p
        break
    }
}

# Temporarily bookmark the start of "Melee Damage" so that the upcoming
# replace (which coalesces adjacent same-formatted runs in its paragraph)
# doesn't also swallow the "Melee Damage" run into the merge.
$full2 = $d.Content.Text
$meleeIdx = $full2.IndexOf("Melee Damage")
$guardSpot = $d.Range($meleeIdx, $meleeIdx)
$d.Bookmarks.Add("ZZGuard", $guardSpot) | Out-Null

# A no-op replace over "Increase Sum" forces the run split (caused by the
# bookmark that used to live there) to heal back into one run.
$meleePara.Range.Find.Execute("Increase Sum", $true, $false, $false, $false, $false, $true, 1, $false, "Increase Sum", 1) | Out-Null

# Remove the temporary guard bookmark now that the merge is done.
$d.Bookmarks("ZZGuard").Delete()
